# Add a new "alimento" user row test data -> update row 4 of the Usuarios
# sheet to reflect the finished "añadir alimento" feature test user.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: rename test user from Lucia/MeInvento/M to Fulano/Fulanito/H,
# and refresh their activity/measurement sample data.
$ws.Range("B4").Value = "Fulano"
$ws.Range("C4").Value = "Fulanito"
$ws.Range("E4").Value = "H"
$ws.Range("F4").Value = 32
$ws.Range("G4").Value = 185
$ws.Range("H4").Value = 83
